$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2397003745318352
$ws.Range("C2").Value2 = 0.4269662921348314
$ws.Range("J2").Value2 = 0.0299625468164794
$ws.Range("P2").Value2 = 0.1760299625468165
$ws.Range("S2").Value2 = 0.1273408239700375
$ws.Range("B3").Value2 = 0.02479338842975207
$ws.Range("C3").Value2 = 0.03305785123966942
$ws.Range("J3").Value2 = 0.04958677685950413
$ws.Range("P3").Value2 = 0.6611570247933884
$ws.Range("S3").Value2 = 0.231404958677686
$ws.Range("J4").Value2 = 0.04545454545454546
$ws.Range("P4").Value2 = 0.6136363636363636
$ws.Range("S4").Value2 = 0.3409090909090909
$ws.Range("P5").Value2 = 1
$ws.Range("B6").Value2 = 0.07511737089201878
$ws.Range("D6").Value2 = 0.009389671361502348
$ws.Range("F6").Value2 = 0.09389671361502347
$ws.Range("J6").Value2 = 0.2910798122065728
$ws.Range("O6").Value2 = 0.0187793427230047
$ws.Range("Q6").Value2 = 0.1408450704225352
$ws.Range("R6").Value2 = 0.05164319248826291
$ws.Range("S6").Value2 = 0.3192488262910798
$ws.Range("B7").Value2 = 0.1228070175438596
$ws.Range("D7").Value2 = 0.005847953216374269
$ws.Range("F7").Value2 = 0.08187134502923976
$ws.Range("J7").Value2 = 0.1111111111111111
$ws.Range("O7").Value2 = 0.005847953216374269
$ws.Range("Q7").Value2 = 0.1286549707602339
$ws.Range("R7").Value2 = 0.0935672514619883
$ws.Range("S7").Value2 = 0.4502923976608187
$ws.Range("B8").Value2 = 0.07808564231738035
$ws.Range("D8").Value2 = 0.02267002518891688
$ws.Range("E8").Value2 = 0.002518891687657431
$ws.Range("F8").Value2 = 0.06801007556675064
$ws.Range("J8").Value2 = 0.1486146095717884
$ws.Range("O8").Value2 = 0.02267002518891688
$ws.Range("Q8").Value2 = 0.1838790931989925
$ws.Range("R8").Value2 = 0.07304785894206549
$ws.Range("S8").Value2 = 0.4005037783375315
$ws.Range("B9").Value2 = 0.07222222222222222
$ws.Range("D9").Value2 = 0.03888888888888889
$ws.Range("F9").Value2 = 0.03888888888888889
$ws.Range("J9").Value2 = 0.1055555555555556
$ws.Range("O9").Value2 = 0.02777777777777778
$ws.Range("Q9").Value2 = 0.1944444444444444
$ws.Range("R9").Value2 = 0.1111111111111111
$ws.Range("S9").Value2 = 0.4111111111111111
$ws.Range("B10").Value2 = 0.1029281277728483
$ws.Range("D10").Value2 = 0.02040816326530612
$ws.Range("E10").Value2 = 0.001774622892635315
$ws.Range("F10").Value2 = 0.06921029281277728
$ws.Range("J10").Value2 = 0.1339840283939663
$ws.Range("O10").Value2 = 0.01774622892635315
$ws.Range("Q10").Value2 = 0.188997338065661
$ws.Range("R10").Value2 = 0.09316770186335403
$ws.Range("S10").Value2 = 0.3717834960070985
$ws.Range("G11").Value2 = 0.1317829457364341
$ws.Range("J11").Value2 = 0.09689922480620156
$ws.Range("K11").Value2 = 0.1821705426356589
$ws.Range("L11").Value2 = 0.5658914728682171
$ws.Range("S11").Value2 = 0.02325581395348837
$ws.Range("G12").Value2 = 0.7058823529411765
$ws.Range("J12").Value2 = 0.2156862745098039
$ws.Range("K12").Value2 = 0.0130718954248366
$ws.Range("L12").Value2 = 0.0392156862745098
$ws.Range("S12").Value2 = 0.0261437908496732
$ws.Range("G13").Value2 = 0.7391304347826086
$ws.Range("J13").Value2 = 0.08695652173913043
$ws.Range("S13").Value2 = 0.1739130434782609
$ws.Range("F15").Value2 = 0.02926829268292683
$ws.Range("H15").Value2 = 0.1560975609756098
$ws.Range("I15").Value2 = 0.1024390243902439
$ws.Range("J15").Value2 = 0.3073170731707317
$ws.Range("K15").Value2 = 0.06829268292682927
$ws.Range("M15").Value2 = 0.01463414634146342
$ws.Range("O15").Value2 = 0.06829268292682927
$ws.Range("S15").Value2 = 0.2536585365853659
$ws.Range("F16").Value2 = 0.01986754966887417
$ws.Range("H16").Value2 = 0.1125827814569536
$ws.Range("I16").Value2 = 0.05298013245033113
$ws.Range("J16").Value2 = 0.4635761589403973
$ws.Range("K16").Value2 = 0.05298013245033113
$ws.Range("M16").Value2 = 0.05298013245033113
$ws.Range("O16").Value2 = 0.03973509933774835
$ws.Range("S16").Value2 = 0.2052980132450331
$ws.Range("F17").Value2 = 0.01851851851851852
$ws.Range("H17").Value2 = 0.1851851851851852
$ws.Range("I17").Value2 = 0.1216931216931217
$ws.Range("J17").Value2 = 0.3756613756613756
$ws.Range("K17").Value2 = 0.1031746031746032
$ws.Range("M17").Value2 = 0.02116402116402116
$ws.Range("O17").Value2 = 0.0582010582010582
$ws.Range("S17").Value2 = 0.1164021164021164
$ws.Range("F18").Value2 = 0.0223463687150838
$ws.Range("H18").Value2 = 0.1843575418994413
$ws.Range("I18").Value2 = 0.111731843575419
$ws.Range("J18").Value2 = 0.4022346368715084
$ws.Range("K18").Value2 = 0.106145251396648
$ws.Range("M18").Value2 = 0.0223463687150838
$ws.Range("O18").Value2 = 0.0446927374301676
$ws.Range("S18").Value2 = 0.106145251396648
$ws.Range("F19").Value2 = 0.02311643835616438
$ws.Range("H19").Value2 = 0.2080479452054795
$ws.Range("I19").Value2 = 0.07448630136986302
$ws.Range("J19").Value2 = 0.348458904109589
$ws.Range("K19").Value2 = 0.1113013698630137
$ws.Range("M19").Value2 = 0.02140410958904109
$ws.Range("O19").Value2 = 0.0761986301369863
$ws.Range("S19").Value2 = 0.136986301369863
